$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '69.001.02'
$ws.Range("E2").Value2 = '  +1.83%  '

$ws.Range("D3").Value2 = '3.818.82'
$ws.Range("E3").Value2 = '  +0.35%  '

$ws.Range("E4").Value2 = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '625.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +4.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '164.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  -0.27%  '

$ws.Range("D7").Value2 = '3.816.32'
$ws.Range("E7").Value2 = '  +0.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.518'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = '  +0.14%  '

$ws.Range("E10").Value2 = '  +1.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.453'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = '  +0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '6.64'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = '  +3.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '0.0000249'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = '  +0.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '35.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = '  -0.03%  '

$ws.Range("D15").Value2 = '4.466.52'
$ws.Range("E15").Value2 = '  +0.45%  '

$ws.Range("D16").Value2 = '3.829.56'
$ws.Range("E16").Value2 = '  +0.87%  '

$ws.Range("D17").Value2 = '69.084.46'
$ws.Range("E17").Value2 = '  +1.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '18.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = '  -1.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '7.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  +1.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '0.113'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = '  +0.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '466.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  +0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '9.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  -1.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '0.708'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = '  +1.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '0.0000152'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  +3.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '84.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = '  +1.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '12.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = '  -0.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '2.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = '  +1.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '10.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = '  +0.37%  '

$ws.Range("E29").Value2 = '  -0.01%  '

$ws.Range("D30").Value2 = '3.977.46'
$ws.Range("E30").Value2 = '  +0.57%  '

$ws.Range("B31").Value2 = 'PancakeSwap'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '2.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = '  -4.16%  '

$ws.Range("B32").Value2 = 'ImmutableX'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '2.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = '  +2.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '7.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = '  -0.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '29.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = '  -0.37%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '9.11'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = '  +0.76%  '

$ws.Range("E36").Value2 = '  -0.13%  '

$ws.Range("E37").Value2 = '  +1.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.148'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = '  +7.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '3.33'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = '  +3.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '5.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  +2.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '0.978'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = '  -1.90%  '

$ws.Range("E42").Value2 = '  +0.08%  '

$ws.Range("E43").Value2 = '  +0.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '155.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = '  +2.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '0.300'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = '  +0.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '1.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = '  +2.68%  '

$ws.Range("B47").Value2 = 'OKB'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '46.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = '  -1.91%  '

$ws.Range("B48").Value2 = 'Arweave'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '42.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = '  -5.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '8.46'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = '  +1.37%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '1.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = '  +1.95%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '381.46'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = '  -3.22%  '

